# Update the "想去人数" (interested-count) figures that changed between
# the previous gh-pages data pull and the new one (commit 456a3b4).
#
# Both the "展览" sheet and the combined "全部类型" sheet carry the same
# rows for this event list, so the same F-column updates are applied to
# each of them.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 522
    "F8"  = 498
    "F9"  = 6580
    "F10" = 178
    "F11" = 142
    "F12" = 1028
    "F13" = 357
    "F14" = 114
    "F15" = 183
    "F16" = 513
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
